$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row below mirrors one coin entry update from the source diff:
# Price (D) is written via NumberFormat "@" so Excel keeps the literal
# text (e.g. "1.000", "29.118.59") instead of coercing it to a number,
# then the style is reset to Normal so no stray cell formatting is left
# behind (matching the original inlineStr-only cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.118.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6962"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.73%  "
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08132"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7270"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.227"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.119.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.784"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("E20").Value = "  -4.53%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.096.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.618"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.986"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("E27").Value = "  -5.39%  "
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.983"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.411"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.491"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.017"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05226"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.022"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.653"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01855"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.679"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9270"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.081.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4269"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.777"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.993.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.199"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.36%  "
